$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.332.21'
$ws.Cells.Item(2, 5).Value = '  +0.29%  '
$ws.Cells.Item(3, 4).Value = '2.008.62'
$ws.Cells.Item(3, 5).Value = '  -1.38%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '258.76'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +4.34%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '0.618'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -2.27%  '
$ws.Cells.Item(7, 5).Value = '  -0.01%  '
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '56.96'
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -6.05%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.383'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -3.25%  '
$ws.Cells.Item(10, 5).Value = '  -4.86%  '
$ws.Cells.Item(11, 5).Value = '  -3.08%  '
$ws.Cells.Item(12, 2).Value = 'Chainlink'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '14.28'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -6.48%  '
$ws.Cells.Item(13, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(13, 4).Value = '2.304.35'
$ws.Cells.Item(13, 5).Value = '  -1.28%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '21.28'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -5.53%  '
$ws.Cells.Item(15, 5).Value = '  -7.60%  '
$ws.Cells.Item(16, 5).Value = '  -5.67%  '
$ws.Cells.Item(17, 4).Value = '2.016.70'
$ws.Cells.Item(17, 5).Value = '  -1.00%  '
$ws.Cells.Item(18, 4).Value = '37.242.73'
$ws.Cells.Item(18, 5).Value = '  +0.18%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '69.77'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.57%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0837'
$ws.Cells.Item(20, 5).Value = '  -3.69%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '231.48'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.10%  '
$ws.Cells.Item(22, 5).Value = '  -3.24%  '
$ws.Cells.Item(23, 5).Value = '  -0.07%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '2.59'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +2.53%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.38%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '164.57'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.30%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '8.95'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -5.78%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '19.57'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -1.75%  '
$ws.Cells.Item(29, 5).Value = '  -6.07%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '1.33'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -4.19%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '0.119'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -2.22%  '
$ws.Cells.Item(32, 5).Value = '  -3.69%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '4.58'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -5.96%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '4.50'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.35%  '
$ws.Cells.Item(35, 5).Value = '  -6.47%  '
$ws.Cells.Item(36, 5).Value = '  +0.43%  '
$ws.Cells.Item(37, 5).Value = '  -0.08%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '3.37'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -2.64%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '5.45'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.15%  '
$ws.Cells.Item(40, 5).Value = '  +3.14%  '
$ws.Cells.Item(41, 5).Value = '  -0.62%  '
$ws.Cells.Item(42, 2).Value = 'Cronos'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.0927'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -5.68%  '
$ws.Cells.Item(43, 2).Value = 'VeChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.0211'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -2.01%  '
$ws.Cells.Item(44, 4).Value = '1.417.54'
$ws.Cells.Item(44, 5).Value = '  +1.85%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '15.73'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -7.19%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '89.46'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -3.52%  '
$ws.Cells.Item(47, 5).Value = '  -4.08%  '
$ws.Cells.Item(48, 5).Value = '  +2.06%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '7.01'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -7.51%  '
$ws.Cells.Item(50, 4).Value = '2.196.32'
$ws.Cells.Item(50, 5).Value = '  -1.35%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '1.95'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -9.35%  '
